$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated forecast error table values for rows 7-11 (Q6-Q10), columns B-G
$ws.Range("B7").Value = -0.1276315789473684
$ws.Range("C7").Value = 0.4939473684210525
$ws.Range("D7").Value = 0.4701500000000001
$ws.Range("E7").Value = 0.6856748500564973
$ws.Range("F7").Value = 0.6827347006035729
$ws.Range("G7").Value = 38

$ws.Range("B8").Value = -0.1083783783783784
$ws.Range("C8").Value = 0.4997297297297297
$ws.Range("D8").Value = 0.4612459459459459
$ws.Range("E8").Value = 0.6791509007179082
$ws.Range("F8").Value = 0.6796956570305465
$ws.Range("G8").Value = 37

$ws.Range("B9").Value = -0.1855
$ws.Range("C9").Value = 0.5124999999999998
$ws.Range("D9").Value = 0.489935
$ws.Range("E9").Value = 0.6999535698887462
$ws.Range("F9").Value = 0.6924591950736918
$ws.Range("G9").Value = 20

$ws.Range("B10").Value = -0.08615384615384616
$ws.Range("C10").Value = 0.5569230769230769
$ws.Range("D10").Value = 0.6175384615384615
$ws.Range("E10").Value = 0.7858361543849083
$ws.Range("F10").Value = 0.8129938341457297
$ws.Range("G10").Value = 13

$ws.Range("B11").Value = -0.516
$ws.Range("C11").Value = 0.532
$ws.Range("D11").Value = 0.5188399999999999
$ws.Range("E11").Value = 0.7203054907468079
$ws.Range("F11").Value = 0.561898567358914
$ws.Range("G11").Value = 5
